$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix F8/F9 formulas to reference column D instead of E
$ws.Range("F8").Formula = "=D8/`$C`$3"
$ws.Range("F9").Formula = "=D9/`$C`$4"

# Update D14 value (Total used after), which ripples into E14/F14
$ws.Range("D14").Value = 13.81

# Turn F13:F14 into a shared formula referencing D13/D14 over C4
$ws.Range("F13:F14").Formula = "=D13/`$C`$4"

# Update the active cell selection on the sheet
$ws.Range("D16").Select() | Out-Null
